# Restore/update the "From" value for rule R30 (row 10) from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("C10").Value = 1
